$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 273 - this shifts existing rows 273:336 down to 274:337
$ws.Rows.Item(273).Insert()

# Populate the newly inserted row 273 with the new weekly record
$ws.Range("A273").Value = 7
$ws.Range("B273").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C273").Value = "Ñuble"
$ws.Range("D273").Value = 45135
$ws.Range("D273").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E273").Value = 16
$ws.Range("F273").Value = 100112032
$ws.Range("G273").Value = "Zapallo italiano"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 80
$ws.Range("K273").Value = 14000
$ws.Range("L273").Value = 14000
$ws.Range("M273").Value = 14000
$ws.Range("N273").Value = "$/caja 50 unidades"
$ws.Range("O273").Value = "Región de Arica y Parinacota"
$ws.Range("P273").Value = 280
$ws.Range("Q273").Value = 50
$ws.Range("R273").Value = "Hortaliza"
